$wb = $excel.ActiveWorkbook

# --- Sheet: Step1_Data ---
$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("AA2").Value = 0.00443113698871727
$ws.Range("AD2").Value = 0.01309382709019463
$ws.Range("D2").Value = 0.02409684053321085
$ws.Range("E2").Value = 0.2066701687279229
$ws.Range("F2").Value = 0.2368748211430806
$ws.Range("G2").Value = 0.06926493272476017
$ws.Range("H2").Value = 0.04940570950799093
$ws.Range("I2").Value = 0.0214843957674714
$ws.Range("K2").Value = 0.01623278281605089
$ws.Range("L2").Value = 0.001869262307145624
$ws.Range("M2").Value = 0.006627632163958432
$ws.Range("N2").Value = 0.0007127135213769319
$ws.Range("O2").Value = 0.03949892543450891
$ws.Range("P2").Value = 0.194259396732218
$ws.Range("R2").Value = 0.09982349075711436
$ws.Range("U2").Value = 0.002239612674201767
$ws.Range("X2").Value = 0.01020288532981355
$ws.Range("Z2").Value = 0.00321146578026286
$ws.Range("D3").Value = 0.06835532538458586
$ws.Range("E3").Value = 0.1280394973148889
$ws.Range("F3").Value = 0.1608763598446709
$ws.Range("G3").Value = 0.02482800652806261
$ws.Range("H3").Value = 0.06666678567991512
$ws.Range("I3").Value = 0.02740556180344403
$ws.Range("K3").Value = 0.0415372253663019
$ws.Range("L3").Value = 0.01554962728700707
$ws.Range("M3").Value = 0.0004366459704207822
$ws.Range("N3").Value = 0.0231587589295313
$ws.Range("O3").Value = 0.06452825899822606
$ws.Range("P3").Value = 0.1888555571318137
$ws.Range("R3").Value = 0.1240539539754357
$ws.Range("S3").Value = 0.0183028499472165
$ws.Range("T3").Value = 0.004514514844829508
$ws.Range("U3").Value = 0.02836154878585952
$ws.Range("X3").Value = 0.01452952220779053
$ws.Range("AC4").Value = 0.01193994933802523
$ws.Range("AD4").Value = 0.01021468890447012
$ws.Range("D4").Value = 0.3700895542072284
$ws.Range("E4").Value = 0.09452707538155031
$ws.Range("F4").Value = 0.1148616482276693
$ws.Range("G4").Value = 0.02853373487350601
$ws.Range("H4").Value = 0.03437254933558269
$ws.Range("J4").Value = 0.03440368922253909
$ws.Range("M4").Value = 0.01118834869787645
$ws.Range("N4").Value = 0.0006453502481165082
$ws.Range("O4").Value = 0.1675113455497269
$ws.Range("P4").Value = 0.01646238075655795
$ws.Range("Q4").Value = 0.07218678808315897
$ws.Range("R4").Value = 0.01282288044175003
$ws.Range("S4").Value = 0.005365167170670121
$ws.Range("X4").Value = 0.001462190466922244
$ws.Range("Z4").Value = 0.01341265909464961
$ws.Range("AA5").Value = 0.005115542882690435
$ws.Range("AE5").Value = 0.000921291549764527
$ws.Range("AF5").Value = 0.002925285466877849
$ws.Range("AG5").Value = 0.009194098108931192
$ws.Range("E5").Value = 0.01824474329318395
$ws.Range("F5").Value = 0.1124307001356686
$ws.Range("G5").Value = 0.2325910761248871
$ws.Range("H5").Value = 0.06415969934756204
$ws.Range("I5").Value = 0.06838840851783386
$ws.Range("J5").Value = 0.0153118610386494
$ws.Range("L5").Value = 0.05579236107804458
$ws.Range("O5").Value = 0.03945348464695239
$ws.Range("P5").Value = 0.07585462781139984
$ws.Range("Q5").Value = 0.1461680220434769
$ws.Range("S5").Value = 0.1324289274599313
$ws.Range("V5").Value = 0.02101987049414588
$ws.Range("AA6").Value = 0.005250472128443546
$ws.Range("AC6").Value = 0.02401112004678332
$ws.Range("AD6").Value = 0.03245450039445103
$ws.Range("D6").Value = 0.02540266555443305
$ws.Range("E6").Value = 0.2427256254655629
$ws.Range("F6").Value = 0.2485568686426603
$ws.Range("G6").Value = 0.04957984733566602
$ws.Range("H6").Value = 0.03630563650211673
$ws.Range("J6").Value = 0.0004509256048816311
$ws.Range("K6").Value = 0.01901688560989294
$ws.Range("M6").Value = 0.000748868331838288
$ws.Range("O6").Value = 0.05181354439524598
$ws.Range("P6").Value = 0.1402060763429704
$ws.Range("R6").Value = 0.1118913547139555
$ws.Range("T6").Value = 0.003703203722462036
$ws.Range("X6").Value = 0.004006117884106966
$ws.Range("Z6").Value = 0.003876287324529264

# --- Sheet: Step2_Sj ---
$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("AA2").Value = 0.9869061729098054
$ws.Range("AB2").Value = 0.9869061729098054
$ws.Range("AC2").Value = 0.9869061729098054
$ws.Range("D2").Value = 0.02409684053321085
$ws.Range("E2").Value = 0.2307670092611337
$ws.Range("F2").Value = 0.4676418304042144
$ws.Range("G2").Value = 0.5369067631289746
$ws.Range("H2").Value = 0.5863124726369655
$ws.Range("I2").Value = 0.6077968684044369
$ws.Range("J2").Value = 0.6077968684044369
$ws.Range("K2").Value = 0.6240296512204878
$ws.Range("L2").Value = 0.6258989135276334
$ws.Range("M2").Value = 0.6325265456915918
$ws.Range("N2").Value = 0.6332392592129688
$ws.Range("O2").Value = 0.6727381846474777
$ws.Range("P2").Value = 0.8669975813796956
$ws.Range("Q2").Value = 0.8669975813796956
$ws.Range("R2").Value = 0.96682107213681
$ws.Range("S2").Value = 0.96682107213681
$ws.Range("T2").Value = 0.96682107213681
$ws.Range("U2").Value = 0.9690606848110118
$ws.Range("V2").Value = 0.9690606848110118
$ws.Range("W2").Value = 0.9690606848110118
$ws.Range("X2").Value = 0.9792635701408253
$ws.Range("Y2").Value = 0.9792635701408253
$ws.Range("Z2").Value = 0.9824750359210882
$ws.Range("D3").Value = 0.06835532538458586
$ws.Range("E3").Value = 0.1963948226994748
$ws.Range("F3").Value = 0.3572711825441456
$ws.Range("G3").Value = 0.3820991890722082
$ws.Range("H3").Value = 0.4487659747521233
$ws.Range("I3").Value = 0.4761715365555674
$ws.Range("J3").Value = 0.4761715365555674
$ws.Range("K3").Value = 0.5177087619218693
$ws.Range("L3").Value = 0.5332583892088764
$ws.Range("M3").Value = 0.5336950351792972
$ws.Range("N3").Value = 0.5568537941088285
$ws.Range("O3").Value = 0.6213820531070546
$ws.Range("P3").Value = 0.8102376102388683
$ws.Range("Q3").Value = 0.8102376102388683
$ws.Range("R3").Value = 0.934291564214304
$ws.Range("S3").Value = 0.9525944141615206
$ws.Range("T3").Value = 0.95710892900635
$ws.Range("U3").Value = 0.9854704777922095
$ws.Range("V3").Value = 0.9854704777922095
$ws.Range("W3").Value = 0.9854704777922095
$ws.Range("AA4").Value = 0.9778453617575044
$ws.Range("AB4").Value = 0.9778453617575044
$ws.Range("AC4").Value = 0.9897853110955296
$ws.Range("AD4").Value = 0.9999999999999997
$ws.Range("AE4").Value = 0.9999999999999997
$ws.Range("AF4").Value = 0.9999999999999997
$ws.Range("AG4").Value = 0.9999999999999997
$ws.Range("AH4").Value = 0.9999999999999997
$ws.Range("AI4").Value = 0.9999999999999997
$ws.Range("AJ4").Value = 0.9999999999999997
$ws.Range("D4").Value = 0.3700895542072284
$ws.Range("E4").Value = 0.4646166295887788
$ws.Range("F4").Value = 0.579478277816448
$ws.Range("G4").Value = 0.608012012689954
$ws.Range("H4").Value = 0.6423845620255366
$ws.Range("I4").Value = 0.6423845620255366
$ws.Range("J4").Value = 0.6767882512480757
$ws.Range("K4").Value = 0.6767882512480757
$ws.Range("L4").Value = 0.6767882512480757
$ws.Range("M4").Value = 0.6879765999459522
$ws.Range("N4").Value = 0.6886219501940687
$ws.Range("O4").Value = 0.8561332957437955
$ws.Range("P4").Value = 0.8725956765003534
$ws.Range("Q4").Value = 0.9447824645835123
$ws.Range("R4").Value = 0.9576053450252624
$ws.Range("S4").Value = 0.9629705121959325
$ws.Range("T4").Value = 0.9629705121959325
$ws.Range("U4").Value = 0.9629705121959325
$ws.Range("V4").Value = 0.9629705121959325
$ws.Range("W4").Value = 0.9629705121959325
$ws.Range("X4").Value = 0.9644327026628547
$ws.Range("Y4").Value = 0.9644327026628547
$ws.Range("Z4").Value = 0.9778453617575044
$ws.Range("AA5").Value = 0.9869593248744262
$ws.Range("AB5").Value = 0.9869593248744262
$ws.Range("AC5").Value = 0.9869593248744262
$ws.Range("AD5").Value = 0.9869593248744262
$ws.Range("AE5").Value = 0.9878806164241907
$ws.Range("AF5").Value = 0.9908059018910685
$ws.Range("AG5").Value = 0.9999999999999998
$ws.Range("AH5").Value = 0.9999999999999998
$ws.Range("AI5").Value = 0.9999999999999998
$ws.Range("AJ5").Value = 0.9999999999999998
$ws.Range("E5").Value = 0.01824474329318395
$ws.Range("F5").Value = 0.1306754434288525
$ws.Range("G5").Value = 0.3632665195537396
$ws.Range("H5").Value = 0.4274262189013017
$ws.Range("I5").Value = 0.4958146274191355
$ws.Range("J5").Value = 0.5111264884577849
$ws.Range("K5").Value = 0.5111264884577849
$ws.Range("L5").Value = 0.5669188495358294
$ws.Range("M5").Value = 0.5669188495358294
$ws.Range("N5").Value = 0.5669188495358294
$ws.Range("O5").Value = 0.6063723341827818
$ws.Range("P5").Value = 0.6822269619941816
$ws.Range("Q5").Value = 0.8283949840376585
$ws.Range("R5").Value = 0.8283949840376585
$ws.Range("S5").Value = 0.9608239114975898
$ws.Range("T5").Value = 0.9608239114975898
$ws.Range("U5").Value = 0.9608239114975898
$ws.Range("V5").Value = 0.9818437819917357
$ws.Range("W5").Value = 0.9818437819917357
$ws.Range("X5").Value = 0.9818437819917357
$ws.Range("Y5").Value = 0.9818437819917357
$ws.Range("Z5").Value = 0.9818437819917357
$ws.Range("AA6").Value = 0.9435343795587654
$ws.Range("AB6").Value = 0.9435343795587654
$ws.Range("AC6").Value = 0.9675454996055487
$ws.Range("AD6").Value = 0.9999999999999998
$ws.Range("AE6").Value = 0.9999999999999998
$ws.Range("AF6").Value = 0.9999999999999998
$ws.Range("AG6").Value = 0.9999999999999998
$ws.Range("AH6").Value = 0.9999999999999998
$ws.Range("AI6").Value = 0.9999999999999998
$ws.Range("AJ6").Value = 0.9999999999999998
$ws.Range("D6").Value = 0.02540266555443305
$ws.Range("E6").Value = 0.268128291019996
$ws.Range("F6").Value = 0.5166851596626563
$ws.Range("G6").Value = 0.5662650069983223
$ws.Range("H6").Value = 0.602570643500439
$ws.Range("I6").Value = 0.602570643500439
$ws.Range("J6").Value = 0.6030215691053206
$ws.Range("K6").Value = 0.6220384547152135
$ws.Range("L6").Value = 0.6220384547152135
$ws.Range("M6").Value = 0.6227873230470518
$ws.Range("N6").Value = 0.6227873230470518
$ws.Range("O6").Value = 0.6746008674422977
$ws.Range("P6").Value = 0.8148069437852681
$ws.Range("Q6").Value = 0.8148069437852681
$ws.Range("R6").Value = 0.9266982984992236
$ws.Range("S6").Value = 0.9266982984992236
$ws.Range("T6").Value = 0.9304015022216856
$ws.Range("U6").Value = 0.9304015022216856
$ws.Range("V6").Value = 0.9304015022216856
$ws.Range("W6").Value = 0.9304015022216856
$ws.Range("X6").Value = 0.9344076201057926
$ws.Range("Y6").Value = 0.9344076201057926
$ws.Range("Z6").Value = 0.9382839074303219

# --- Sheet: Step3_DataPts_0.5 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("D2").Value = 6
$ws.Range("F2").Value = 0.5369067631289746
$ws.Range("G2").Value = 4
$ws.Range("D3").Value = 10
$ws.Range("F3").Value = 0.5177087619218693
$ws.Range("G3").Value = 9
$ws.Range("D4").Value = 5
$ws.Range("F4").Value = 0.579478277816448
$ws.Range("G4").Value = 4
$ws.Range("D5").Value = 9
$ws.Range("F5").Value = 0.5111264884577849
$ws.Range("G5").Value = 6
$ws.Range("F6").Value = 0.5166851596626563

# --- Sheet: Step3_DataPts_0.7 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F2").Value = 0.8669975813796956
$ws.Range("F3").Value = 0.8102376102388683
$ws.Range("D4").Value = 14
$ws.Range("F4").Value = 0.8561332957437955
$ws.Range("G4").Value = 13
$ws.Range("F5").Value = 0.8283949840376585
$ws.Range("D6").Value = 15
$ws.Range("F6").Value = 0.8148069437852681
$ws.Range("G6").Value = 13

# --- Sheet: Step3_DataPts_0.8 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F2").Value = 0.8669975813796956
$ws.Range("F3").Value = 0.8102376102388683
$ws.Range("F4").Value = 0.8561332957437955
$ws.Range("F5").Value = 0.8283949840376585
$ws.Range("F6").Value = 0.8148069437852681

# --- Sheet: Step3_DataPts_0.9 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F2").Value = 0.96682107213681
$ws.Range("F3").Value = 0.934291564214304
$ws.Range("D4").Value = 16
$ws.Range("F4").Value = 0.9447824645835123
$ws.Range("G4").Value = 15
$ws.Range("F5").Value = 0.9608239114975898
$ws.Range("F6").Value = 0.9266982984992236
